# Update "想去人数" (want-to-go count) figures across the workbook's sheets
# to match the refreshed data pull (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14084
$ws1.Range("F8").Value = 35
$ws1.Range("F9").Value = 71
$ws1.Range("F10").Value = 796
$ws1.Range("F11").Value = 2176
$ws1.Range("F24").Value = 140

# 演出 (Performance) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 1667
$ws2.Range("F15").Value = 1804

# 本地生活 (Local life) sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 131

# 全部类型 (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 14084
$ws4.Range("F9").Value = 35
$ws4.Range("F10").Value = 71
$ws4.Range("F11").Value = 796
$ws4.Range("F14").Value = 2176
$ws4.Range("F15").Value = 131
$ws4.Range("F34").Value = 1667
$ws4.Range("F39").Value = 140
$ws4.Range("F48").Value = 1804
